$wb = $excel.ActiveWorkbook

# 1. Remove the three "Service Fee" rows (C/D columns, rows 19-21) from the
#    "Combined Price List 25-26 Sum" sheet. This shifts everything below up
#    by 3 rows (the trailing placeholder row moves from 39 to 36).
$priceList = $wb.Worksheets.Item("Combined Price List 25-26 Sum")
$priceList.Rows("19:21").Delete()

# 2. Remove the "Surcharge Order" column (column D) from the
#    "Config 25-26 Sum" sheet.
$config = $wb.Worksheets.Item("Config 25-26 Sum")
$config.Columns("D").Delete()

# 3. Delete the entire "Discount 25-26 Sum" worksheet.
$discount = $wb.Worksheets.Item("Discount 25-26 Sum")
$discount.Delete()
